$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  "+---------------------+-------+-----------------------------------------+--------+",
  "| Timestamp           | UID   | Name                                    | Price  |",
  "+---------------------+-------+-----------------------------------------+--------+",
  "| 2023-01-22 12:16:59 | P001F | APPLE iPhone 14 (Blue, 128 GB)          |  72999 |",
  "| 2023-01-22 12:16:59 | P002F | Xbox Series S Console 512 GB  (White)   |  30490 |",
  "| 2023-01-22 12:18:14 | P001A | Apple iPhone 14 Pro Max 1TB Space Black | 182999 |",
  "| 2023-01-22 12:18:17 | P002A | Apple iPhone 13 (128GB) - Pink          |  61999 |",
  "| 2023-01-22 12:18:18 | P001F | APPLE iPhone 14 (Blue, 128 GB)          |  72999 |",
  "| 2023-01-22 12:18:19 | P002F | Xbox Series S Console 512 GB  (White)   |  30490 |",
  "| 2023-01-22 12:40:35 | P001A | Apple iPhone 14 Pro Max 1TB Space Black | 182999 |",
  "| 2023-01-22 12:40:38 | P002A | Apple iPhone 13 (128GB) - Pink          |  61999 |",
  "| 2023-01-22 12:40:39 | P001F | APPLE iPhone 14 (Blue, 128 GB)          |  72999 |",
  "| 2023-01-22 12:40:39 | P002F | Xbox Series S Console 512 GB  (White)   |  30490 |",
  "| 2023-01-22 12:45:09 | P001A | Apple iPhone 14 Pro Max 1TB Space Black | 182999 |",
  "| 2023-01-22 12:45:12 | P002A | Apple iPhone 13 (128GB) - Pink          |  61999 |",
  "| 2023-01-22 12:45:13 | P001F | APPLE iPhone 14 (Blue, 128 GB)          |  72999 |",
  "| 2023-01-22 12:45:13 | P002F | Xbox Series S Console 512 GB  (White)   |  30490 |",
  "| 2023-01-22 12:48:26 | P001A | Apple iPhone 14 Pro Max 1TB Space Black | 182999 |",
  "| 2023-01-22 12:48:29 | P002A | Apple iPhone 13 (128GB) - Pink          |  61999 |",
  "| 2023-01-22 12:48:30 | P001F | APPLE iPhone 14 (Blue, 128 GB)          |  72999 |",
  "| 2023-01-22 12:48:31 | P002F | Xbox Series S Console 512 GB  (White)   |  30490 |",
  "| 2023-01-22 12:54:45 | P001A | Apple iPhone 14 Pro Max 1TB Space Black | 182999 |",
  "| 2023-01-22 12:54:48 | P002A | Apple iPhone 13 (128GB) - Pink          |  61999 |",
  "| 2023-01-22 12:54:48 | P001F | APPLE iPhone 14 (Blue, 128 GB)          |  72999 |",
  "| 2023-01-22 12:54:49 | P002F | Xbox Series S Console 512 GB  (White)   |  30490 |",
  "| 2023-01-22 13:00:56 | P002A | Apple iPhone 13 (128GB) - Pink          |  61999 |",
  "| 2023-01-22 13:00:57 | P001F | APPLE iPhone 14 (Blue, 128 GB)          |  72999 |",
  "| 2023-01-22 13:00:58 | P002F | Xbox Series S Console 512 GB  (White)   |  30490 |",
  "+---------------------+-------+-----------------------------------------+--------+"
)

$startRow = 31
for ($i = 0; $i -lt $rows.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 1).Value = $rows[$i]
}

$endRow = $startRow + $rows.Length - 1
$ws.Range("A$startRow`:A$endRow").Select()
$excel.ActiveWindow.ScrollRow = 5
